# Parcial II: PRograma Parcial
# Subo programa con parcial, punto 2 parece que ya, falta verificar punto 1
#
# Adds a new "tiempo" / "temp/concentra" label row below the Runge-Kutta
# table on sheet "x" (row 14), and moves the active selection to A17.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("x")
$ws.Activate()

# New row 14: time / temp-concentration labels (used as future chart axis
# labels for the RK4 table above).
$ws.Range("E14").Value = "tiempo"
$ws.Range("F14").Value = "temp/concentra"
$ws.Range("K14").Value = "temp/concentra"

# Move the selection like the author left it before committing.
$ws.Range("A17").Select()
